$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting rows 99:211 down to 100:212
$ws.Rows("99:99").Insert()

# Fill the new row 99 with data (copy the static columns from row 100, which is the
# row that used to be row 99 before the shift, to keep A,B,C,E,F,G,N,O,Q,R identical)
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44638
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112045
$ws.Range("G99").Value = "Zapallo"
$ws.Range("H99").Value = "Camote"
$ws.Range("I99").Value = "1a (cosecha)"
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 300
$ws.Range("L99").Value = 300
$ws.Range("M99").Value = 300
$ws.Range("N99").Value = "$/kilo (volumen en unidades)"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 300
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"
